$d = $word.ActiveDocument

function New-PackageXml($bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $bodyXml + '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
}

# ------------------------------------------------------------------
# 1) Turn the first two (empty) paragraphs into:
#      - a paragraph holding two single-space runs
#      - a paragraph (keeping the original pPr/rFonts) that now hosts
#        the "_GoBack" bookmark
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)

# Merge paragraph 1 and paragraph 2 into a single paragraph mark so that
# InsertXML can cleanly replace the whole span with exactly two new
# paragraphs (InsertXML tends to leave a stray empty paragraph behind
# when asked to replace a range spanning more than one existing
# paragraph mark).
$mergeRange = $d.Range($p1.Range.End - 1, $p1.Range.End)
$mergeRange.Delete()

$p1 = $d.Paragraphs.Item(1)
$introXml = New-PackageXml (
    '<w:p>' +
        '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p>' +
    '<w:p>' +
        '<w:pPr><w:rPr><w:rFonts w:ascii="微软雅黑" w:eastAsia="微软雅黑" w:hAnsi="微软雅黑"/></w:rPr></w:pPr>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>'
)
$p1.Range.InsertXML($introXml)

# ------------------------------------------------------------------
# 2) Remove the old "_GoBack" bookmark that used to sit in the middle
#    of the "praseInt(...)" run, and merge the two runs it used to
#    separate ("...retur" + "n") back into a single run.
# ------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("praseInt")) {
        $target = $p
        break
    }
}

$mergedParaXml = New-PackageXml (
    '<w:p>' +
        '<w:pPr><w:ind w:firstLineChars="50" w:firstLine="105"/><w:rPr><w:rFonts w:ascii="微软雅黑" w:eastAsia="微软雅黑" w:hAnsi="微软雅黑"/><w:szCs w:val="21"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:rFonts w:ascii="微软雅黑" w:eastAsia="微软雅黑" w:hAnsi="微软雅黑" w:hint="eastAsia"/><w:szCs w:val="21"/></w:rPr><w:t>praseInt（2，1） return</w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:ascii="微软雅黑" w:eastAsia="微软雅黑" w:hAnsi="微软雅黑"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> N</w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:ascii="微软雅黑" w:eastAsia="微软雅黑" w:hAnsi="微软雅黑" w:hint="eastAsia"/><w:szCs w:val="21"/></w:rPr><w:t>a</w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:ascii="微软雅黑" w:eastAsia="微软雅黑" w:hAnsi="微软雅黑"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">N </w:t></w:r>' +
    '</w:p>'
)
$target.Range.InsertXML($mergedParaXml)
